$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 1000
$ws.Cells.Item(70, 9).Value = 1000
$ws.Cells.Item(70, 11).Value = 3000
$ws.Cells.Item(70, 13).Value = -2730

$ws.Cells.Item(73, 8).Value = 1000
$ws.Cells.Item(73, 9).Value = 1000
$ws.Cells.Item(73, 11).Value = 3000
$ws.Cells.Item(73, 13).Value = -2064

$ws.Cells.Item(98, 8).Value = 1659.125
$ws.Cells.Item(98, 9).Value = 1512.1666
$ws.Cells.Item(98, 10).Value = 2100
$ws.Cells.Item(98, 11).Value = 1512.1666
$ws.Cells.Item(98, 12).Value = 2100
$ws.Cells.Item(98, 13).Value = -14.16660000000002
$ws.Cells.Item(98, 14).Value = -5096

$ws.Cells.Item(122, 8).Value = 1659.125
$ws.Cells.Item(122, 9).Value = 1512.1666
$ws.Cells.Item(122, 10).Value = 2100
$ws.Cells.Item(122, 11).Value = 4536.4998
$ws.Cells.Item(122, 12).Value = 6300
$ws.Cells.Item(122, 13).Value = -2086.4998
$ws.Cells.Item(122, 14).Value = -11200

$ws.Cells.Item(132, 8).Value = 2133050.2
$ws.Cells.Item(132, 9).Value = 2882.8948
$ws.Cells.Item(132, 10).Value = 12251345
$ws.Cells.Item(132, 11).Value = 8648.6844
$ws.Cells.Item(132, 12).Value = 36754035
$ws.Cells.Item(132, 13).Value = -6118.6844
$ws.Cells.Item(132, 14).Value = -36759095

$ws.Cells.Item(137, 8).Value = 2274836
$ws.Cells.Item(137, 9).Value = 3335077
$ws.Cells.Item(137, 10).Value = 2891.0715
$ws.Cells.Item(137, 11).Value = 10005231
$ws.Cells.Item(137, 12).Value = 8673.2145
$ws.Cells.Item(137, 13).Value = -10002681
$ws.Cells.Item(137, 14).Value = -13773.2145

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(12, 8).Value = 6000
$ws.Cells.Item(12, 9).Value = 5000
$ws.Cells.Item(12, 10).Value = 7000
$ws.Cells.Item(12, 11).Value = 5000
$ws.Cells.Item(12, 12).Value = 7000
$ws.Cells.Item(12, 13).Value = -4827
$ws.Cells.Item(12, 14).Value = -7346

$ws.Cells.Item(14, 8).Value = 6750
$ws.Cells.Item(14, 9).Value = 6750
$ws.Cells.Item(14, 11).Value = 6750
$ws.Cells.Item(14, 13).Value = -6575

$ws.Cells.Item(16, 8).Value = 11003.5
$ws.Cells.Item(16, 9).Value = 10000
$ws.Cells.Item(16, 10).Value = 12007
$ws.Cells.Item(16, 11).Value = 10000
$ws.Cells.Item(16, 12).Value = 12007
$ws.Cells.Item(16, 13).Value = -9713
$ws.Cells.Item(16, 14).Value = -12581

$ws.Cells.Item(27, 8).Value = 3000
$ws.Cells.Item(27, 9).Value = 3000
$ws.Cells.Item(27, 11).Value = 3000
$ws.Cells.Item(27, 13).Value = -2816

$ws.Cells.Item(32, 8).Value = 20748.715
$ws.Cells.Item(32, 9).Value = 22215.72
$ws.Cells.Item(32, 11).Value = 22215.72
$ws.Cells.Item(32, 13).Value = -21928.72

$ws.Cells.Item(122, 8).Value = 1866.4
$ws.Cells.Item(122, 9).Value = 1583
$ws.Cells.Item(122, 10).Value = 3000
$ws.Cells.Item(122, 11).Value = 4749
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = -2299
$ws.Cells.Item(122, 14).Value = -13900

$ws.Cells.Item(132, 8).Value = 38889.5
$ws.Cells.Item(132, 9).Value = 31032.559
$ws.Cells.Item(132, 10).Value = 52246.3
$ws.Cells.Item(132, 11).Value = 93097.677
$ws.Cells.Item(132, 12).Value = 156738.9
$ws.Cells.Item(132, 13).Value = -90567.677
$ws.Cells.Item(132, 14).Value = -161798.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2141.0952
$ws.Cells.Item(107, 9).Value = 1493.875
$ws.Cells.Item(107, 10).Value = 2539.3845
$ws.Cells.Item(107, 11).Value = 1493.875
$ws.Cells.Item(107, 12).Value = 2539.3845
$ws.Cells.Item(107, 13).Value = 426.125
$ws.Cells.Item(107, 14).Value = -6379.3845

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 2795.4546
$ws.Cells.Item(62, 10).Value = 2600
$ws.Cells.Item(62, 12).Value = 2600
$ws.Cells.Item(62, 14).Value = -3848

$ws.Cells.Item(65, 8).Value = 2795.4546
$ws.Cells.Item(65, 10).Value = 2600
$ws.Cells.Item(65, 12).Value = 13000
$ws.Cells.Item(65, 14).Value = -19240

$ws.Cells.Item(69, 8).Value = 7122.1816
$ws.Cells.Item(69, 9).Value = 5816
$ws.Cells.Item(69, 10).Value = 13000
$ws.Cells.Item(69, 11).Value = 5816
$ws.Cells.Item(69, 12).Value = 13000
$ws.Cells.Item(69, 13).Value = -5067
$ws.Cells.Item(69, 14).Value = -14498

$ws.Cells.Item(72, 8).Value = 7122.1816
$ws.Cells.Item(72, 9).Value = 5816
$ws.Cells.Item(72, 10).Value = 13000
$ws.Cells.Item(72, 11).Value = 17448
$ws.Cells.Item(72, 12).Value = 39000
$ws.Cells.Item(72, 13).Value = -13704
$ws.Cells.Item(72, 14).Value = -46488

$ws.Cells.Item(99, 8).Value = 3111.361
$ws.Cells.Item(99, 9).Value = 2673.423
$ws.Cells.Item(99, 10).Value = 4250
$ws.Cells.Item(99, 11).Value = 2673.423
$ws.Cells.Item(99, 12).Value = 4250
$ws.Cells.Item(99, 13).Value = -1175.423
$ws.Cells.Item(99, 14).Value = -7246

$ws.Cells.Item(122, 8).Value = 1754.5883
$ws.Cells.Item(122, 9).Value = 1455.8823
$ws.Cells.Item(122, 10).Value = 2053.2942
$ws.Cells.Item(122, 11).Value = 4367.6469
$ws.Cells.Item(122, 12).Value = 6159.882599999999
$ws.Cells.Item(122, 13).Value = -1917.6469
$ws.Cells.Item(122, 14).Value = -11059.8826

$ws.Cells.Item(126, 8).Value = 3111.361
$ws.Cells.Item(126, 9).Value = 2673.423
$ws.Cells.Item(126, 10).Value = 4250
$ws.Cells.Item(126, 11).Value = 8020.268999999999
$ws.Cells.Item(126, 12).Value = 12750
$ws.Cells.Item(126, 13).Value = -5550.268999999999
$ws.Cells.Item(126, 14).Value = -17690

$ws.Cells.Item(132, 8).Value = 30613.314
$ws.Cells.Item(132, 9).Value = 2233.8125
$ws.Cells.Item(132, 10).Value = 54511.844
$ws.Cells.Item(132, 11).Value = 6701.4375
$ws.Cells.Item(132, 12).Value = 163535.532
$ws.Cells.Item(132, 13).Value = -4171.4375
$ws.Cells.Item(132, 14).Value = -168595.532

$ws.Cells.Item(134, 8).Value = 53578.047
$ws.Cells.Item(134, 9).Value = 2516.0715
$ws.Cells.Item(134, 10).Value = 155702
$ws.Cells.Item(134, 11).Value = 7548.2145
$ws.Cells.Item(134, 12).Value = 467106
$ws.Cells.Item(134, 13).Value = -5013.2145
$ws.Cells.Item(134, 14).Value = -472176

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(115, 8).Value = 2816.8572

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21, 8).Value = 4000
$ws.Cells.Item(21, 9).Value = 3000
$ws.Cells.Item(21, 10).Value = 5000
$ws.Cells.Item(21, 11).Value = 3000
$ws.Cells.Item(21, 12).Value = 5000
$ws.Cells.Item(21, 13).Value = -2827
$ws.Cells.Item(21, 14).Value = -5346

$ws.Cells.Item(30, 8).Value = 4000
$ws.Cells.Item(30, 9).Value = 3000
$ws.Cells.Item(30, 10).Value = 5000
$ws.Cells.Item(30, 11).Value = 3000
$ws.Cells.Item(30, 12).Value = 5000
$ws.Cells.Item(30, 13).Value = -2895
$ws.Cells.Item(30, 14).Value = -5210

$ws.Cells.Item(123, 8).Value = 20755.875
$ws.Cells.Item(123, 10).Value = 20755.875
$ws.Cells.Item(123, 12).Value = 20755.875
$ws.Cells.Item(123, 14).Value = -25655.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 13).ClearContents() | Out-Null

$ws.Cells.Item(61, 8).Value = 4475
$ws.Cells.Item(61, 9).Value = 4475
$ws.Cells.Item(61, 11).Value = 4475
$ws.Cells.Item(61, 13).Value = -4273

$ws.Cells.Item(113, 8).Value = 4475
$ws.Cells.Item(113, 9).Value = 4475
$ws.Cells.Item(113, 11).Value = 4475
$ws.Cells.Item(113, 13).Value = -2305

$ws.Cells.Item(132, 8).Value = 87148.336
$ws.Cells.Item(132, 9).Value = 3060
$ws.Cells.Item(132, 10).Value = 147211.42
$ws.Cells.Item(132, 11).Value = 9180
$ws.Cells.Item(132, 12).Value = 441634.26
$ws.Cells.Item(132, 13).Value = -6650
$ws.Cells.Item(132, 14).Value = -446694.26

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 14).ClearContents() | Out-Null

$ws.Cells.Item(114, 8).Value = 0
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 14).ClearContents() | Out-Null

$ws.Cells.Item(126, 8).Value = 1340.8334
$ws.Cells.Item(126, 9).Value = 979.5238000000001
$ws.Cells.Item(126, 10).Value = 2183.889
$ws.Cells.Item(126, 11).Value = 2938.5714
$ws.Cells.Item(126, 12).Value = 6551.667
$ws.Cells.Item(126, 13).Value = -468.5714000000003
$ws.Cells.Item(126, 14).Value = -11491.667
